# Apply the diff to price_report_template.xlsx:
#   - cell B6 text "[d.month]" -> "[d.start_date] [EN DASH] [d.end_date]"
#   - selection on sheet changes from A8 to B6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$enDash = [char]0x2013
$ws.Range("B6").Value = "[d.start_date] " + $enDash + " [d.end_date]"

$ws.Range("B6").Select()
